# 6.4.2.1.xlsx — add the "2020" column (J) mirroring the existing
# "2019" column (I): same per-row formatting, new figures, and the
# J7 "by type of source" subtotal is a formula (Total - groundwater).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone column I's cell formatting (font/borders/number format/alignment)
# into column J for every data row, so the new column visually matches
# the existing year columns.
for ($r = 4; $r -le 18; $r++) {
  $ws.Range("I$r").Copy()
  $ws.Range("J$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Header year
$ws.Range("J4").Value = 2020

# Total freshwater withdrawal
$ws.Range("J5").Value = 8017.9

# "By type of source" section
$ws.Range("J6").ClearContents()
$ws.Range("J8").Value = 249.8
$ws.Range("J7").Formula = "=J5-J8"

# "By territory" section
$ws.Range("J9").ClearContents()
$ws.Range("J10").Value = 757.6
$ws.Range("J11").Value = 984.4
$ws.Range("J12").Value = 646.20000000000005
$ws.Range("J13").Value = 667.6
$ws.Range("J14").Value = 1147
$ws.Range("J15").Value = 961.1
$ws.Range("J16").Value = 2664.5
$ws.Range("J17").Value = 132.5
$ws.Range("J18").Value = 57

# Selection moves to J19, matching the saved view state.
[void]$ws.Range("J19").Select()
